# Update build timestamp strings across the workbook for the
# "mines - January 30" release (new build: February 02 2026 12.49.33 EST)

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 (built on January 30 2026 16.19.47 EST)"
$newStamp = "January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: mines - " + $newStamp

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"' + "Global Energy Monitor, Coal mine boundaries and methane sources for Sima Coal Mine, China, M1203, version '" + "mines - " + $newStamp + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $wsData.Range("S" + $row).Value = "mines - " + $newStamp
}
